$wb = $excel.ActiveWorkbook

# --- Sheet: FBS ---
$wsFBS = $wb.Worksheets.Item("FBS")
$wsFBS.Range("Z5").Value = -106
$wsFBS.Range("P8").Value = 5.2
$wsFBS.Range("O8").Value = 50.33000000000001
$wsFBS.Range("N8").Value = "NW"
$wsFBS.Range("Q8").Value = "NW"
$wsFBS.Range("U8").Value = -5.4
$wsFBS.Range("Z8").Value = -114
$wsFBS.Range("Y8").Value = 62.5
$wsFBS.Range("AE8").Value = 0.01626016260162602
$wsFBS.Range("Q9").Value = "SE"
$wsFBS.Range("P9").Value = 2.4
$wsFBS.Range("M9").Value = "NNW"
$wsFBS.Range("N9").Value = "SE"
$wsFBS.Range("O9").Value = 39.56
$wsFBS.Range("U9").Value = -6.9
$wsFBS.Range("Y9").Value = 47.5
$wsFBS.Range("AE9").Value = 0
$wsFBS.Range("M10").Value = "SSE"
$wsFBS.Range("P10").Value = 2
$wsFBS.Range("N10").Value = "E"
$wsFBS.Range("O10").Value = 39.26
$wsFBS.Range("Q10").Value = "E"
$wsFBS.Range("U10").Value = -3.1
$wsFBS.Range("Q11").Value = "ESE"
$wsFBS.Range("O11").Value = 53.24
$wsFBS.Range("P11").Value = 6.2
$wsFBS.Range("M11").Value = "SE"
$wsFBS.Range("N11").Value = "ESE"
$wsFBS.Range("U11").Value = -2
$wsFBS.Range("N12").Value = "SE"
$wsFBS.Range("M12").Value = "SE"
$wsFBS.Range("Q12").Value = "SE"
$wsFBS.Range("Z12").Value = -108
$wsFBS.Range("Y12").Value = 51.5
$wsFBS.Range("AE12").Value = 0
$wsFBS.Range("M13").Value = "WNW"
$wsFBS.Range("N13").Value = "WNW"
$wsFBS.Range("Q13").Value = "WNW"
$wsFBS.Range("N14").Value = "SE"
$wsFBS.Range("M14").Value = "SE"
$wsFBS.Range("Q14").Value = "SE"
$wsFBS.Range("Q16").Value = "NW"
$wsFBS.Range("M19").Value = "S"
$wsFBS.Range("N19").Value = "S"
$wsFBS.Range("Q19").Value = "S"
$wsFBS.Range("N20").Value = "SE"
$wsFBS.Range("M20").Value = "SE"
$wsFBS.Range("Q20").Value = "SE"
$wsFBS.Range("Z20").Value = -106
$wsFBS.Range("N21").Value = "SE"
$wsFBS.Range("Q21").Value = "SE"
$wsFBS.Range("N22").Value = "WNW"
$wsFBS.Range("Q22").Value = "NW"
$wsFBS.Range("Q24").Value = "N"
$wsFBS.Range("N25").Value = "SE"
$wsFBS.Range("Q25").Value = "SE"
$wsFBS.Range("AF26").Value = 4
$wsFBS.Range("AB26").Value = -19
$wsFBS.Range("N30").Value = "S"
$wsFBS.Range("M30").Value = "S"
$wsFBS.Range("Q30").Value = "S"
$wsFBS.Range("AB30").Value = -8
$wsFBS.Range("AF30").Value = 0.5
$wsFBS.Range("M32").Value = "SE"
$wsFBS.Range("N32").Value = "SE"
$wsFBS.Range("Q32").Value = "SE"
$wsFBS.Range("AB33").Value = -9.5
$wsFBS.Range("AF33").Value = 0
$wsFBS.Range("Z33").Value = -110
$wsFBS.Range("N34").Value = "SE"
$wsFBS.Range("M34").Value = "SE"
$wsFBS.Range("Q34").Value = "SE"
$wsFBS.Range("M36").Value = "SE"
$wsFBS.Range("N36").Value = "SE"
$wsFBS.Range("Q36").Value = "SE"
$wsFBS.Range("N37").Value = "NW"
$wsFBS.Range("Q37").Value = "NW"
$wsFBS.Range("N38").Value = "NW"
$wsFBS.Range("Q38").Value = "NW"
$wsFBS.Range("M39").Value = "SE"
$wsFBS.Range("Z40").Value = -106
$wsFBS.Range("M41").Value = "NW"
$wsFBS.Range("N41").Value = "NW"
$wsFBS.Range("Q41").Value = "NW"
$wsFBS.Range("M42").Value = "S"
$wsFBS.Range("Q42").Value = "S"
$wsFBS.Range("N46").Value = "NW"
$wsFBS.Range("Q46").Value = "NW"
$wsFBS.Range("N47").Value = "WNW"
$wsFBS.Range("N49").Value = "WNW"
$wsFBS.Range("Q49").Value = "SW"
$wsFBS.Range("Q50").Value = "NW"
$wsFBS.Range("M51").Value = "WNW"
$wsFBS.Range("N51").Value = "WNW"
$wsFBS.Range("Q51").Value = "WNW"

# --- Sheet: Other ---
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("P5").Value = "SE"
$wsOther.Range("S5").Value = "SE"
$wsOther.Range("P10").Value = "SE"
$wsOther.Range("O10").Value = "SE"
$wsOther.Range("S10").Value = "SE"
$wsOther.Range("O12").Value = "SE"
$wsOther.Range("P12").Value = "SE"
$wsOther.Range("S12").Value = "SE"
$wsOther.Range("S15").Value = "NE"
$wsOther.Range("Q16").Value = 49.88
$wsOther.Range("O16").Value = "S"
$wsOther.Range("R16").Value = 4
$wsOther.Range("S16").Value = "S"
$wsOther.Range("P18").Value = "SE"
$wsOther.Range("O18").Value = "SE"
$wsOther.Range("S18").Value = "SE"
$wsOther.Range("P24").Value = "SE"
$wsOther.Range("O24").Value = "SE"
$wsOther.Range("S24").Value = "SE"
$wsOther.Range("S26").Value = "NW"
$wsOther.Range("P29").Value = "SE"
$wsOther.Range("P33").Value = "WNW"
$wsOther.Range("O33").Value = "WNW"
$wsOther.Range("S33").Value = "WNW"
$wsOther.Range("P38").Value = "NW"
$wsOther.Range("O38").Value = "NW"
$wsOther.Range("S38").Value = "NW"
$wsOther.Range("P39").Value = "NW"
$wsOther.Range("O39").Value = "NW"
$wsOther.Range("S39").Value = "NW"
$wsOther.Range("P43").Value = "NW"
$wsOther.Range("S43").Value = "NW"
$wsOther.Range("O44").Value = "S"
$wsOther.Range("S44").Value = "S"
$wsOther.Range("P49").Value = "S"
$wsOther.Range("S49").Value = "S"
$wsOther.Range("P50").Value = "SE"
$wsOther.Range("O50").Value = "SE"
$wsOther.Range("S50").Value = "SE"
$wsOther.Range("O26").Value = "WNW"
$wsOther.Range("P26").Value = "NW"

# --- Timestamp shared string (all AK2:AK51 on FBS reference the same text) ---
$wsFBS.Range("AK2:AK51").Value = "2024-11-14T10:01:56.607984"
